# Apply the commit: correct emails, remove hyperlinks, add new row (shan),
# and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the email values -------------------------------------------------
# Row 2 (ABOBAKAR) already had the correct address typo-fixed to
# ranaabobakar777@gmail.com
$ws.Range("C2").Value = "ranaabobakar777@gmail.com"

# Row 3 (SHAHZEB) previously had 19014156-022@uog.edu.pk -> now the same
# corrected address
$ws.Range("C3").Value = "ranaabobakar777@gmail.com"

# Row 4 (ZEESHAN) already used ranaabobakar777@gmail.com - keep as is
$ws.Range("C4").Value = "ranaabobakar777@gmail.com"

# --- Remove all hyperlinks on the sheet (and their relationships) --------
$ws.Hyperlinks.Delete() | Out-Null

# --- Add the new row 5 (shan) ---------------------------------------------
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "shan"
$ws.Range("C5").Value = "shan@gmail.com"
$ws.Range("D5").Value = "WAS"

# --- Move the active selection to C7 (matches the saved selection state) -
$ws.Range("C7").Select() | Out-Null
